$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 90, shifting existing rows 90-125 down to 91-126
$ws.Range("A90").EntireRow.Insert()

# Populate the newly inserted row 90 with the new data record
$ws.Range("A90").Value = 10
$ws.Range("B90").Value = "Vega Modelo de Temuco"
$ws.Range("C90").Value = "La Araucanía"
$ws.Range("D90").Value = 44809
$ws.Range("E90").Value = 9
$ws.Range("F90").Value = 100112035
$ws.Range("G90").Value = "Bruselas (repollito)"
$ws.Range("H90").Value = "Sin especificar"
$ws.Range("I90").Value = "Primera"
$ws.Range("J90").Value = 160
$ws.Range("K90").Value = 24000
$ws.Range("L90").Value = 25000
$ws.Range("M90").Value = 24500
$ws.Range("N90").Value = '$/malla 10 kilos'
$ws.Range("O90").Value = "Provincia de Quillota"
$ws.Range("P90").Value = 2450
$ws.Range("Q90").Value = 10
$ws.Range("R90").Value = "Hortaliza"
